$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = "ah ma è ronco "
$ws.Range("B27").Value = "Matteo Zanlucchi | SBARX"
$ws.Range("C27").Value = "Carlo  Stedile | Mai una gioia"
$ws.Range("D27").Value = "Nicolo  Speziali | FC GORILLAZ"
$ws.Range("E27").Value = "Matteo Simoncelli | IMONTAGNA"
$ws.Range("F27").Value = "Francesco Cristoforetti | Vigili del Fusto"
